$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmtShort = "[$-14009]d\ mmmm\ yyyy;@"

# New header cell (J1) - introduces shared string index 12 "Studentoffer"
$ws.Cells.Item(1, 10).Value = "Studentoffer"

# Row 3
$ws.Cells.Item(3, 1).Value = 8015332963
$ws.Cells.Item(3, 1).Font.Color = 16777215
$ws.Cells.Item(3, 2).Value = 8234
$ws.Cells.Item(3, 3).Value = "Chennai"
$ws.Cells.Item(3, 4).Value = "Mumbai"
$ws.Cells.Item(3, 5).Value = 45936
$ws.Cells.Item(3, 5).NumberFormat = $dateFmtShort
$ws.Cells.Item(3, 6).Value = 45942
$ws.Cells.Item(3, 6).NumberFormat = $dateFmtShort
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = "Economy"
$ws.Cells.Item(3, 10).Value = "Student"

# Row 4
$ws.Cells.Item(4, 1).Value = 8015332963
$ws.Cells.Item(4, 1).Font.Color = 16777215
$ws.Cells.Item(4, 2).Value = 8234
$ws.Cells.Item(4, 3).Value = "Coimbatore"
$ws.Cells.Item(4, 4).Value = "Hyderabad"
$ws.Cells.Item(4, 5).Value = 45931
$ws.Cells.Item(4, 5).NumberFormat = $dateFmtShort
$ws.Cells.Item(4, 6).Value = 45961
$ws.Cells.Item(4, 6).NumberFormat = $dateFmtShort
$ws.Cells.Item(4, 7).Value = 2
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = "Premium Economy"
$ws.Cells.Item(4, 10).Value = "Senior Citizen"

# Row 5
$ws.Cells.Item(5, 1).Value = 8015332963
$ws.Cells.Item(5, 1).Font.Color = 16777215
$ws.Cells.Item(5, 2).Value = 8234
$ws.Cells.Item(5, 3).Value = "Punjab"
$ws.Cells.Item(5, 4).Value = "Tiruchirapalli"
$ws.Cells.Item(5, 5).Value = 46013
$ws.Cells.Item(5, 5).NumberFormat = $dateFmtShort
$ws.Cells.Item(5, 6).Value = 46040
$ws.Cells.Item(5, 6).NumberFormat = $dateFmtShort
$ws.Cells.Item(5, 7).Value = 1
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = "Business"
$ws.Cells.Item(5, 10).Value = "Armed Forces"

# Column widths for the newly introduced / resized columns
$ws.Columns.Item(9).ColumnWidth = 15.343333333333334
$ws.Columns.Item(10).ColumnWidth = 13.676666666666666
$ws.Columns.Item(11).ColumnWidth = 19.17666666666667
$ws.Columns.Item(12).ColumnWidth = 27.01

$ws.Range("J6").Select() | Out-Null
